$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.853.66'
$ws.Range('E2').Value = '  -1.21%  '
$ws.Range('D3').Value = '2.230.81'
$ws.Range('E3').Value = '  -1.79%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.627'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '76.05'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.89%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -2.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.87'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0941'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.23%  '
$ws.Range('E12').Value = '  -2.36%  '
$ws.Range('E13').Value = '  -1.28%  '
$ws.Range('D14').Value = '2.566.56'
$ws.Range('E14').Value = '  -1.95%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.78'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E16').Value = '  -2.07%  '
$ws.Range('D17').Value = '2.230.15'
$ws.Range('E17').Value = '  -2.15%  '
$ws.Range('D18').Value = '41.778.91'
$ws.Range('E18').Value = '  -1.34%  '
$ws.Range('D19').Value = '0.0₃0974'
$ws.Range('E19').Value = '  -2.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.11'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.10'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.14%  '
$ws.Range('E22').Value = '  -0.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '229.78'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.10%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('E25').Value = '  -5.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.07'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.28%  '
$ws.Range('E27').Value = '  -4.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.26'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +14.11%  '
$ws.Range('E29').Value = '  -1.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '169.35'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.42'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0851'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '33.05'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.119'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.125'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('E36').Value = '  -2.18%  '
$ws.Range('E37').Value = '  +1.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0295'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.46'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.18'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.84'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '111.78'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +15.02%  '
$ws.Range('E43').Value = '  -5.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '59.71'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.61%  '
$ws.Range('E45').Value = '  -4.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1000'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.997'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.13'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.98%  '
$ws.Range('E49').Value = '  -1.40%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.20'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -14.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.26'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.64%  '
